# The edit performs a cyclic rotation of the species-observation data that
# lives in rows 2, 3 and 4 of the "Artfynd" sheet:
#   Row 2 receives the data that used to be in Row 3
#   Row 3 receives the data that used to be in Row 4
#   Row 4 receives the data that used to be in Row 2
# Columns C, D, I, P, S..AY (administrative/location/report columns) stay put.
# Only columns A, B, E, F, G, H, K, L, M, N, Q, R travel with the rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the "before" state of the three rows for the columns that move ---
function Get-RowData($row) {
    $data = New-Object System.Collections.Hashtable
    $data["A"] = $ws.Range("A$row").Value2
    $data["B"] = $ws.Range("B$row").Value2
    $data["E"] = $ws.Range("E$row").Value2
    $data["F"] = $ws.Range("F$row").Value2
    $data["G"] = $ws.Range("G$row").Value2
    $data["H"] = $ws.Range("H$row").Value2
    $data["K"] = $ws.Range("K$row").Value2
    $data["L"] = $ws.Range("L$row").Value2
    $data["M"] = $ws.Range("M$row").Value2
    $data["N"] = $ws.Range("N$row").Value2
    $data["Q"] = $ws.Range("Q$row").Value2
    $data["R"] = $ws.Range("R$row").Value2
    return $data
}

# Writing "" through .Value simply deletes/clears the cell in this COM
# layer, which would lose the "present but blank" inline-string cells that
# the original file had (e.g. K3/L3/N3). Using the classic Excel
# "force text" apostrophe prefix makes Excel store a genuine empty-string
# text cell instead of clearing it, matching the source file's shape.
function Set-CellValue($range, $value) {
    if ($value -eq $null) {
        $range.ClearContents()
    } elseif ($value -eq "") {
        $range.Value = "'"
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

function Set-RowData($row, $data) {
    Set-CellValue $ws.Range("A$row") $data["A"]
    Set-CellValue $ws.Range("B$row") $data["B"]
    Set-CellValue $ws.Range("E$row") $data["E"]
    Set-CellValue $ws.Range("F$row") $data["F"]
    Set-CellValue $ws.Range("G$row") $data["G"]
    Set-CellValue $ws.Range("H$row") $data["H"]
    Set-CellValue $ws.Range("Q$row") $data["Q"]
    Set-CellValue $ws.Range("R$row") $data["R"]

    Set-CellValue $ws.Range("K$row") $data["K"]
    Set-CellValue $ws.Range("L$row") $data["L"]
    Set-CellValue $ws.Range("M$row") $data["M"]
    Set-CellValue $ws.Range("N$row") $data["N"]
}

$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row4 = Get-RowData 4

Set-RowData 2 $row3
Set-RowData 3 $row4
Set-RowData 4 $row2
